$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.094.44"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.596.10"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'311.88"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").Value = "'99.40"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.598"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.585"
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").Value = "'39.14"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0843"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").Value = "'54.03"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "'8.16"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "2.993.35"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "2.595.61"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").Value = "'14.87"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").Value = "46.243.36"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").Value = "'6.76"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").Value = "'12.82"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").Value = "'291.47"
$ws.Range("E23").Value = "  +14.97%  "
$ws.Range("D24").Value = "'73.09"
$ws.Range("E24").Value = "  +3.22%  "
$ws.Range("D25").Value = "'3.06"
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("D27").Value = "'29.36"
$ws.Range("E27").Value = "  +6.24%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "'4.05"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").Value = "'10.88"
$ws.Range("E30").Value = "  +4.96%  "
$ws.Range("D31").Value = "'39.36"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").Value = "'6.25"
$ws.Range("E33").Value = "  +2.52%  "
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").Value = "'155.92"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("D36").Value = "'0.0838"
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("D38").Value = "'2.79"
$ws.Range("E38").Value = "  -4.87%  "
$ws.Range("E39").Value = "  +4.76%  "
$ws.Range("D40").Value = "'0.124"
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").Value = "'15.75"
$ws.Range("E41").Value = "  +1.86%  "
$ws.Range("E42").Value = "  +3.81%  "
$ws.Range("D43").Value = "'3.60"
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("D44").Value = "'4.01"
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("D45").Value = "'20.96"
$ws.Range("E45").Value = "  +10.16%  "
$ws.Range("D46").Value = "2.108.50"
$ws.Range("E46").Value = "  +3.25%  "
$ws.Range("D47").Value = "'98.02"
$ws.Range("E47").Value = "  +8.22%  "
$ws.Range("D48").Value = "'0.998"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "'9.50"
$ws.Range("E49").Value = "  +4.36%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.202"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'108.87"
$ws.Range("E51").Value = "  +0.12%  "
